{"js": "// Update the 25 two-digit-by-two-digit multiplication prompts in the\n// practice-sheet table. Each old expression is unique in the document,\n// so a direct search+replace on body text is unambiguous and order\n// independent.\nconst replacements = [\n  [\"35\u00d760=\", \"87\u00d716=\"],\n  [\"15\u00d734=\", \"99\u00d776=\"],\n  [\"66\u00d729=\", \"80\u00d771=\"],\n  [\"36\u00d786=\", \"36\u00d758=\"],\n  [\"94\u00d758=\", \"96\u00d771=\"],\n  [\"97\u00d718=\", \"60\u00d733=\"],\n  [\"50\u00d750=\", \"33\u00d797=\"],\n  [\"28\u00d751=\", \"66\u00d775=\"],\n  [\"66\u00d758=\", \"28\u00d714=\"],\n  [\"53\u00d750=\", \"31\u00d744=\"],\n  [\"44\u00d723=\", \"63\u00d737=\"],\n  [\"26\u00d745=\", \"72\u00d714=\"],\n  [\"37\u00d730=\", \"55\u00d746=\"],\n  [\"95\u00d719=\", \"26\u00d764=\"],\n  [\"37\u00d714=\", \"75\u00d799=\"],\n  [\"94\u00d762=\", \"72\u00d724=\"],\n  [\"33\u00d777=\", \"79\u00d725=\"],\n  [\"19\u00d785=\", \"69\u00d755=\"],\n  [\"30\u00d712=\", \"47\u00d717=\"],\n  [\"33\u00d727=\", \"91\u00d748=\"],\n  [\"26\u00d731=\", \"16\u00d729=\"],\n  [\"20\u00d770=\", \"22\u00d769=\"],\n  [\"59\u00d778=\", \"25\u00d719=\"],\n  [\"33\u00d717=\", \"57\u00d785=\"],\n  [\"92\u00d719=\", \"41\u00d747=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit-by-two-digit multiplication prompts in the\n# practice-sheet table. Each old expression is unique in the document,\n# so a Find/Replace pass over the whole document body is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"35\u00d760=\"; new=\"87\u00d716=\"},\n    @{old=\"15\u00d734=\"; new=\"99\u00d776=\"},\n    @{old=\"66\u00d729=\"; new=\"80\u00d771=\"},\n    @{old=\"36\u00d786=\"; new=\"36\u00d758=\"},\n    @{old=\"94\u00d758=\"; new=\"96\u00d771=\"},\n    @{old=\"97\u00d718=\"; new=\"60\u00d733=\"},\n    @{old=\"50\u00d750=\"; new=\"33\u00d797=\"},\n    @{old=\"28\u00d751=\"; new=\"66\u00d775=\"},\n    @{old=\"66\u00d758=\"; new=\"28\u00d714=\"},\n    @{old=\"53\u00d750=\"; new=\"31\u00d744=\"},\n    @{old=\"44\u00d723=\"; new=\"63\u00d737=\"},\n    @{old=\"26\u00d745=\"; new=\"72\u00d714=\"},\n    @{old=\"37\u00d730=\"; new=\"55\u00d746=\"},\n    @{old=\"95\u00d719=\"; new=\"26\u00d764=\"},\n    @{old=\"37\u00d714=\"; new=\"75\u00d799=\"},\n    @{old=\"94\u00d762=\"; new=\"72\u00d724=\"},\n    @{old=\"33\u00d777=\"; new=\"79\u00d725=\"},\n    @{old=\"19\u00d785=\"; new=\"69\u00d755=\"},\n    @{old=\"30\u00d712=\"; new=\"47\u00d717=\"},\n    @{old=\"33\u00d727=\"; new=\"91\u00d748=\"},\n    @{old=\"26\u00d731=\"; new=\"16\u00d729=\"},\n    @{old=\"20\u00d770=\"; new=\"22\u00d769=\"},\n    @{old=\"59\u00d778=\"; new=\"25\u00d719=\"},\n    @{old=\"33\u00d717=\"; new=\"57\u00d785=\"},\n    @{old=\"92\u00d719=\"; new=\"41\u00d747=\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
